# Update Name of Algo
# Apply updated numeric values to the RandomForest result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -6.247699999999996
$ws.Range("E6").Value = 12.4768
$ws.Range("E7").Value = 11.98889999999999
$ws.Range("D8").Value = -8.945299999999987
$ws.Range("E8").Value = 12.5697
$ws.Range("C12").Value = -14.68050000000001
$ws.Range("D12").Value = -8.219600000000002
$ws.Range("D14").Value = -8.794699999999999
$ws.Range("E19").Value = 12.71799999999999
$ws.Range("E21").Value = 12.57219999999999
$ws.Range("D22").Value = -8.148399999999992
$ws.Range("E24").Value = 12.78539999999999
